# Refresh the "Price" (D) / "Volume(1h)" (E) columns with the latest
# coinranking.com snapshot, and fix the MXToken/BTSEToken row ordering
# (rows 8-9 swapped back to the correct rank).
#
# Price/Volume cells hold numeric-looking text (e.g. "308.33", "-1.97%")
# that must stay literal text, not be parsed into a Number/Percentage -
# so those assignments use a leading apostrophe, same as typing '308.33
# into Excel, to force a text entry instead of numeric coercion.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "'308.33"
$ws.Range("E2").Value = "'-1.97%"
$ws.Range("D3").Value = "'40.99"
$ws.Range("E3").Value = "'-0.43%"
$ws.Range("D4").Value = "'5.040"
$ws.Range("E4").Value = "'-1.87%"
$ws.Range("D5").Value = "'0.07640"
$ws.Range("E5").Value = "'-3.15%"
$ws.Range("D6").Value = "'4.237"
$ws.Range("E6").Value = "'-2.49%"
$ws.Range("D7").Value = "'1.617"
$ws.Range("E7").Value = "'-3.20%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9096"
$ws.Range("E8").Value = "'-1.17%"
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").Value = "'2.452"
$ws.Range("E9").Value = "'-4.66%"
$ws.Range("D10").Value = "'0.1008"
$ws.Range("E10").Value = "'-9.00%"
$ws.Range("D11").Value = "'0.1767"
$ws.Range("E11").Value = "'-2.13%"
$ws.Range("D12").Value = "'0.09038"
$ws.Range("E12").Value = "'-0.76%"
$ws.Range("D13").Value = "'0.04402"
$ws.Range("E13").Value = "'-0.99%"
$ws.Range("D14").Value = "'0.1053"
$ws.Range("E14").Value = "'-0.71%"
$ws.Range("D15").Value = "'0.001252"
$ws.Range("E15").Value = "'-1.31%"
$ws.Range("D16").Value = "'0.005795"
$ws.Range("E16").Value = "'-1.54%"
$ws.Range("D17").Value = "'3.365"
$ws.Range("E17").Value = "'0.37%"
$ws.Range("D19").Value = "'6.777"
$ws.Range("E19").Value = "'-6.77%"
$ws.Range("D20").Value = "'0.1356"
$ws.Range("E20").Value = "'-2.09%"
$ws.Range("D21").Value = "'0.2721"
$ws.Range("E21").Value = "'2.36%"
$ws.Range("D22").Value = "'0.04150"
$ws.Range("E22").Value = "'-0.67%"
$ws.Range("D23").Value = "'0.001213"
$ws.Range("E23").Value = "'-3.18%"
$ws.Range("D24").Value = "'0.004093"
$ws.Range("E24").Value = "'-1.55%"
$ws.Range("D25").Value = "'0.0001299"
$ws.Range("E25").Value = "'5.48%"
$ws.Range("D26").Value = "'0.0003005"
$ws.Range("E26").Value = "'-0.26%"
$ws.Range("D38").Value = "'0.02414"
$ws.Range("E38").Value = "'-2.10%"
$ws.Range("D39").Value = "'0.05175"
$ws.Range("E39").Value = "'-3.21%"
$ws.Range("E40").Value = "'-3.74%"
$ws.Range("E41").Value = "'-3.57%"
$ws.Range("D42").Value = "'0.007087"
$ws.Range("E42").Value = "'-6.42%"
$ws.Range("D43").Value = "'0.001947"
$ws.Range("E43").Value = "'-5.75%"
$ws.Range("E44").Value = "'-3.01%"
$ws.Range("D45").Value = "'0.3058"
$ws.Range("E45").Value = "'-2.21%"
$ws.Range("D46").Value = "'0.00006380"
$ws.Range("E46").Value = "'-6.41%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-1.19%"
$ws.Range("D48").Value = "'0.004397"
$ws.Range("E48").Value = "'6.04%"
$ws.Range("D49").Value = "'0.006124"
$ws.Range("E49").Value = "'78.69%"
$ws.Range("D50").Value = "'0.00002099"
$ws.Range("E50").Value = "'-1.19%"
$ws.Range("D51").Value = "'0.0001999"
$ws.Range("E51").Value = "'-1.19%"
